$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on D2:D51 so numeric-looking price strings
# (e.g. "0.9993", "240.53") are stored as text, matching the source data,
# instead of being auto-converted to numbers by the usual Excel input rules.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.377.43"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.846.44"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "240.53"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "0.6276"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.07496"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "24.47"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.845.92"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "0.6809"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "0.00001056"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "6.184"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "29.413.07"
$ws.Range("D19").Value = "229.11"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "7.484"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "159.48"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "8.422"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "17.54"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "0.06459"
$ws.Range("E28").Value = "  +15.36%  "
$ws.Range("D29").Value = "1.414"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("D31").Value = "4.096"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "4.096"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "1.832"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").Value = "0.6983"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "2.584"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "1.267.67"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").Value = "2.841"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "0.01829"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").Value = "6.779"
$ws.Range("E40").Value = "  +6.18%  "
$ws.Range("D41").Value = "0.9098"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "2.008.14"
$ws.Range("E43").Value = "  -18.39%  "
$ws.Range("D44").Value = "101.38"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "1.738"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.072"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1166"
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.004"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.3968"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05699"
$ws.Range("E51").Value = "  -0.09%  "

# Clear the temporary text-format override so the cell style matches
# the original (unstyled) data cells.
$ws.Range("D2:D51").ClearFormats()

